# Update the MIAPE - Chromatography template:
#  - bump the template Version from 1.0.0 to 1.0.1
#  - normalize casing of the "Characteristic [Sample description]" column
#    header to "Characteristic [sample description]"
#  - repoint the sample-description ontology reference from
#    MIAPPE:0079 to DPBO:0000175 (Term Source REF / Term Accession Number)

$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump Version value (row "Version" / value cell B4) ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.1"

# --- MIAPE-CC sheet: update the annotation table's header row ---
$wsTable = $wb.Worksheets.Item("MIAPE-CC")
$wsTable.Range("B1").Value = "Characteristic [sample description]"
$wsTable.Range("C1").Value = "Term Source REF (DPBO:0000175)"
$wsTable.Range("D1").Value = "Term Accession Number (DPBO:0000175)"
